$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add preparation material for session 05 (row 6)
$ws.Range("D6").Value = "prep/p04.html"

# Add materials for session 03 (row 4)
$ws.Range("F4").Value = "exercises/e03.html"
$ws.Range("E4").Value = "slides/slides.html#/sitzung-03-warum-wir-unterschiedliche-inhalte-in-manchen-situationen-unterhaltsam-finden"

# Update the active cell selection
$ws.Range("E3").Select()
